$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-26 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-27 Sunday", 2) | Out-Null
$d.Content.Find.Execute("661÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "181÷8=", 2) | Out-Null
$d.Content.Find.Execute("143÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "518÷7=", 2) | Out-Null
$d.Content.Find.Execute("774÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "196÷2=", 2) | Out-Null
$d.Content.Find.Execute("444÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "985÷8=", 2) | Out-Null
$d.Content.Find.Execute("310÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "557÷6=", 2) | Out-Null
$d.Content.Find.Execute("148÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "112÷9=", 2) | Out-Null
$d.Content.Find.Execute("579÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "126÷3=", 2) | Out-Null
$d.Content.Find.Execute("799÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "706÷7=", 2) | Out-Null
$d.Content.Find.Execute("294÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "653÷9=", 2) | Out-Null
$d.Content.Find.Execute("965÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "884÷9=", 2) | Out-Null
$d.Content.Find.Execute("993÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "174÷6=", 2) | Out-Null
$d.Content.Find.Execute("579÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "837÷9=", 2) | Out-Null
$d.Content.Find.Execute("142÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "985÷4=", 2) | Out-Null
$d.Content.Find.Execute("707÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "111÷3=", 2) | Out-Null
$d.Content.Find.Execute("707÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "292÷4=", 2) | Out-Null
$d.Content.Find.Execute("300÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "250÷4=", 2) | Out-Null
$d.Content.Find.Execute("121÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "724÷4=", 2) | Out-Null
$d.Content.Find.Execute("514÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "505÷7=", 2) | Out-Null
$d.Content.Find.Execute("953÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "174÷2=", 2) | Out-Null
$d.Content.Find.Execute("690÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "782÷8=", 2) | Out-Null
$d.Content.Find.Execute("848÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "628÷8=", 2) | Out-Null
$d.Content.Find.Execute("944÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "426÷3=", 2) | Out-Null
$d.Content.Find.Execute("189÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "532÷2=", 2) | Out-Null
$d.Content.Find.Execute("623÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "841÷2=", 2) | Out-Null
$d.Content.Find.Execute("831÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "358÷4=", 2) | Out-Null
